# Set attendance marker cells (0 -> 1) on Sheet1 per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> list of columns that should be set to 1
$changes = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("D", "E")
    6  = @("H")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("H")
    11 = @("G", "H")
    12 = @("H")
    13 = @("D", "E")
    14 = @("H")
    15 = @("H")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
